# DOM and Banner author ids added
# Rows 9-13 were re-matched to the correct OpenAlex records, shifting which
# publication each row describes; M3 (cited_by_count) was refreshed 9 -> 11.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# M3: cited_by_count 9 -> 11 (stored as text, matching the sheet's existing format)
$ws.Cells.Item(3, 13).Formula = "'11"

# Row 9
$ws.Cells.Item(9, 1).Value = "Neusha Barakati, Rocio Zapata Bustos, Dawn K. Coletta, Paul Langlais, Lindsay N. Kohler, Moulun Luo, Janet L. Funk, Wayne T. Willis, Lawrence J. Mandarino"
$ws.Cells.Item(9, 2).Value = "Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona; Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona; Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona; Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona; Center for Disparities in Diabetes, Obesity, and Metabolism, University of Arizona, Health Sciences, Tucson, Arizona; Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona; Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona; Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona; Division of Endocrinology, Department of Medicine, The University of Arizona, Tucson, Arizona"
$ws.Cells.Item(9, 3).Value = "https://openalex.org/W4281290394"
$ws.Cells.Item(9, 4).Value = "Acetylation of Adenine Nucleotide Translocase, Fuel Selection, and Metabolic Flexibility in Human Skeletal Muscle"
$ws.Cells.Item(9, 5).Formula = "'2022-05-09"
$ws.Cells.Item(9, 6).Value = "medRxiv (Cold Spring Harbor Laboratory)"
$ws.Cells.Item(9, 7).Value = "Cold Spring Harbor Laboratory"
$ws.Cells.Item(9, 8).Value = "https://doi.org/10.1101/2022.05.05.22274505"
$ws.Cells.Item(9, 9).Value = "N/A"
$ws.Cells.Item(9, 10).Value = "submittedVersion"
$ws.Cells.Item(9, 11).Value = "green"
$ws.Cells.Item(9, 12).Value = "en"
$ws.Cells.Item(9, 13).Formula = "'0"
$ws.Cells.Item(9, 14).Formula = "'2022"
$ws.Cells.Item(9, 15).Value = "NA"
$ws.Cells.Item(9, 16).Value = "https://doi.org/10.1101/2022.05.05.22274505"
$ws.Cells.Item(9, 17).Value = "article"

# Row 10
$ws.Cells.Item(10, 1).Value = "Jin Ou, Eric M. Lewandowski, Yanmei Hu, Austin A. Lipinski, Ryan T. Morgan, Lian Jacobs, Xiujun Zhang, Melissa J. Bikowitz, Paul Langlais, Haozhou Tan, Jun Wang, Yu. M. Chumakov, John S. Choy"
$ws.Cells.Item(10, 2).Value = " The Catholic University of America;;  University of South Florida;; Department of Medicinal Chemistry, Ernest Mario School of Pharmacy, Rutgers, the State University of New Jersey, Piscataway, NJ, 08854, United States;  University of Arizona;  University of South Florida;; Department of Molecular Medicine, Morsani College of Medicine, University of South Florida, Tampa, FL, 3361, United States; Drug Discovery Department, Moffit Cancer Center, Tampa, FL 33612, United States;  University of South Florida;;  University of South Florida;;  University of Arizona; Department of Medicinal Chemistry, Ernest Mario School of Pharmacy, Rutgers, the State University of New Jersey, Piscataway, NJ, 08854, United States; Department of Medicinal Chemistry, Ernest Mario School of Pharmacy, Rutgers, the State University of New Jersey, Piscataway, NJ, 08854, United States;  University of South Florida;;  The Catholic University of America;"
$ws.Cells.Item(10, 3).Value = "https://openalex.org/W4290988437"
$ws.Cells.Item(10, 4).Value = "A yeast-based system to study SARS-CoV-2 M<sup>pro</sup> structure and to identify nirmatrelvir resistant mutations"
$ws.Cells.Item(10, 5).Formula = "'2022-08-08"
$ws.Cells.Item(10, 6).Value = "bioRxiv (Cold Spring Harbor Laboratory)"
$ws.Cells.Item(10, 7).Value = "Cold Spring Harbor Laboratory"
$ws.Cells.Item(10, 8).Value = "https://doi.org/10.1101/2022.08.06.503039"
$ws.Cells.Item(10, 9).Value = "cc-by-nc"
$ws.Cells.Item(10, 10).Value = "submittedVersion"
$ws.Cells.Item(10, 11).Value = "green"
$ws.Cells.Item(10, 12).Value = "en"
$ws.Cells.Item(10, 13).Formula = "'0"
$ws.Cells.Item(10, 14).Formula = "'2022"
$ws.Cells.Item(10, 15).Value = "https://pubmed.ncbi.nlm.nih.gov/35982672"
$ws.Cells.Item(10, 16).Value = "https://doi.org/10.1101/2022.08.06.503039"
$ws.Cells.Item(10, 17).Value = "article"

# Row 11
$ws.Cells.Item(11, 1).Value = "Velia S. Vizcarra, Kara R. Barber, Gabriela Franca-Solomon, Lisa Majuta, Angela Smith, Paul Langlais, Tally M. Largent‐Milnes, Todd W. Vanderah, Arthur C. Riegel"
$ws.Cells.Item(11, 2).Value = "Department of Pharmacology, College of Medicine, University of Arizona, Tucson, AZ 85721, USA; Department of Pharmacology, College of Medicine, University of Arizona, Tucson, AZ 85721, USA; Department of Pharmacology, College of Medicine, University of Arizona, Tucson, AZ 85721, USA; Department of Pharmacology, College of Medicine, University of Arizona, Tucson, AZ 85721, USA; Department of Pharmacology, College of Medicine, University of Arizona, Tucson, AZ 85721, USA; Department of Medicine, Division of Endocrinology, University of Arizona College of Medicine, Tucson, AZ 85721, USA; Neuroscience Graduate Interdisciplinary Program, University of Arizona, Tucson, AZ 85721, USA; Neuroscience Graduate Interdisciplinary Program, University of Arizona, Tucson, AZ 85721, USA; Department of Pharmacology, College of Medicine, University of Arizona, Tucson, AZ 85721, USA; James C. Wyant College of Optical Sciences, The University of Arizona, Tucson, AZ 85721, USA"
$ws.Cells.Item(11, 3).Value = "https://openalex.org/W4294550774"
$ws.Cells.Item(11, 4).Value = "Targeting 5-HT2A receptors and Kv7 channels in PFC to attenuate chronic neuropathic pain in rats using a spared nerve injury model"
$ws.Cells.Item(11, 5).Formula = "'2022-10-01"
$ws.Cells.Item(11, 6).Value = "Neuroscience Letters"
$ws.Cells.Item(11, 7).Value = "Elsevier BV"
$ws.Cells.Item(11, 8).Value = "https://doi.org/10.1016/j.neulet.2022.136864"
$ws.Cells.Item(11, 9).Value = "cc-by-nc-nd"
$ws.Cells.Item(11, 10).Value = "publishedVersion"
$ws.Cells.Item(11, 11).Value = "hybrid"
$ws.Cells.Item(11, 12).Value = "en"
$ws.Cells.Item(11, 13).Formula = "'0"
$ws.Cells.Item(11, 14).Formula = "'2022"
$ws.Cells.Item(11, 15).Value = "https://pubmed.ncbi.nlm.nih.gov/36063980"
$ws.Cells.Item(11, 16).Value = "https://doi.org/10.1016/j.neulet.2022.136864"
$ws.Cells.Item(11, 17).Value = "article"

# Row 12
$ws.Cells.Item(12, 1).Value = "Sanjay Kumar, Aaron Ramonett, Tasmia Ahmed, Eun-A Kwak, Paola Cruz Flores, Hannah R. Ortiz, Paul Langlais, Karthikeyan Mythreye, Nam Y. Lee"
$ws.Cells.Item(12, 2).Value = "Division of Biology, Indian Institute of Science Education and Research, Tirupati 57507, India; Department of Pharmacology, University of Arizona, Tucson, AZ 8574, USA; Department of Chemistry & Biochemistry, University of Arizona, Tucson, AZ 85724, USA; Department of Pharmacology, University of Arizona, Tucson, AZ 8574, USA; Department of Chemistry & Biochemistry, University of Arizona, Tucson, AZ 85724, USA; Department of Pharmacology, University of Arizona, Tucson, AZ 8574, USA; Department of Medicine, University of Arizona, Tucson, AZ 8572, USA; Department of Pathology, University of Alabama at Birmingham, 3294, USA; Department of Chemistry & Biochemistry, University of Arizona, Tucson, AZ 85724, USA"
$ws.Cells.Item(12, 3).Value = "https://openalex.org/W4205184871"
$ws.Cells.Item(12, 4).Value = "Identification of Mfn2-S249 as a Phosphoregulatory Switch of Mitochondrial Fusion Dynamics"
$ws.Cells.Item(12, 5).Formula = "'2022-01-12"
$ws.Cells.Item(12, 6).Value = "bioRxiv (Cold Spring Harbor Laboratory)"
$ws.Cells.Item(12, 7).Value = "Cold Spring Harbor Laboratory"
$ws.Cells.Item(12, 8).Value = "https://doi.org/10.1101/2022.01.11.475884"
$ws.Cells.Item(12, 9).Value = "N/A"
$ws.Cells.Item(12, 10).Value = "submittedVersion"
$ws.Cells.Item(12, 11).Value = "green"
$ws.Cells.Item(12, 12).Value = "en"
$ws.Cells.Item(12, 13).Formula = "'0"
$ws.Cells.Item(12, 14).Formula = "'2022"
$ws.Cells.Item(12, 15).Value = "NA"
$ws.Cells.Item(12, 16).Value = "https://doi.org/10.1101/2022.01.11.475884"
$ws.Cells.Item(12, 17).Value = "article"

# Row 13
$ws.Cells.Item(13, 1).Value = "Natasha R Cornejo, Bismark Amofah, Austin A. Lipinski, Paul Langlais, Indraneel Ghosh, John C. Jewett"
$ws.Cells.Item(13, 2).Value = "; ; ; ; ; "
$ws.Cells.Item(13, 3).Value = "https://openalex.org/W4283714303"
$ws.Cells.Item(13, 4).Value = "Correction to “Direct Intracellular Delivery of Benzene Diazonium Ions As Observed by Increased Tyrosine Phosphorylation”"
$ws.Cells.Item(13, 5).Formula = "'2022-06-30"
$ws.Cells.Item(13, 6).Value = "Biochemistry"
$ws.Cells.Item(13, 7).Value = "American Chemical Society"
$ws.Cells.Item(13, 8).Value = "https://doi.org/10.1021/acs.biochem.2c00316"
$ws.Cells.Item(13, 9).Value = "N/A"
$ws.Cells.Item(13, 10).Value = "publishedVersion"
$ws.Cells.Item(13, 11).Value = "bronze"
$ws.Cells.Item(13, 12).Value = "en"
$ws.Cells.Item(13, 13).Formula = "'0"
$ws.Cells.Item(13, 14).Formula = "'2022"
$ws.Cells.Item(13, 15).Value = "https://pubmed.ncbi.nlm.nih.gov/35772027"
$ws.Cells.Item(13, 16).Value = "https://doi.org/10.1021/acs.biochem.2c00316"
$ws.Cells.Item(13, 17).Value = "article"
